# Update cryptocurrency price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.738.35'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '1.599.93'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.66'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("D12").Value = '1.824.23'
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").Value = '1.599.82'
$ws.Range("E13").Value = '  -0.04%  '
$ws.Range("E14").Value = '  +1.07%  '
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.04'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").Value = '0.0₃0741'
$ws.Range("E17").Value = '  -1.29%  '
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '208.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("E20").Value = '  +2.39%  '
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.14'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.36'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("E29").Value = '  -2.00%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("E31").Value = '  +0.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.97'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.60%  '
$ws.Range("D33").Value = '1.280.11'
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("E34").Value = '  +1.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.22'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +16.51%  '
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.592'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.73%  '
$ws.Range("E38").Value = '  -1.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.827'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("E42").Value = '  -0.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.35%  '
$ws.Range("D44").Value = '1.736.42'
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.55%  '
$ws.Range("E46").Value = '  -0.18%  '
$ws.Range("E47").Value = '  +2.12%  '
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.07%  '
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("E51").Value = '  +1.60%  '
